# Update the "IS PW - meetings" document:
#  - Append the year "/21" (or "/22") to the existing meeting dates.
#  - Give the "28/12" row an explicit (at-least) row height.
#  - Append a new meeting row for 05/01/22 documenting the class-diagram update.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Fix up the existing date cells (column 1) -------------------------
$dateFixes = @{
    2 = "30/11/21"
    3 = "06/12/21"
    4 = "11/12/21"
    5 = "13/12/21"
    6 = "28/12/21"
}

foreach ($rowIndex in $dateFixes.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $dateFixes[$rowIndex]
}

# --- 2. Give row 6 (28/12/21) an explicit "at least" row height -----------
$row6 = $t.Rows.Item(6)
$row6.HeightRule = 1
$row6.Height = 22.4

# --- 3. Append a new row for the 05/01/22 meeting --------------------------
$t.Rows.Add() | Out-Null
$newRowIndex = $t.Rows.Count

$newRow = $t.Rows.Item($newRowIndex)
$newRow.HeightRule = 1
$newRow.Height = 22.4

$t.Cell($newRowIndex, 1).Range.Text = "05/01/22"
$t.Cell($newRowIndex, 2).Range.Text = "Modifica class diagram."

Write-Output "done"
